$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "MOP_DEF" in column F (new last column), using the same
# formatting as the existing header cells (bold, centered, bordered - style E1).
$ws.Range("F1").Value = "MOP_DEF"

$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats - copy formatting only, keep the value we just set
